# Updates in response to issue #39
# - Add a new "PE_stream" row (StreamImpactItem) to both the "info" sheet
#   and the "GWP" sheet, mirroring the existing "PE" ImpactItem values via
#   formulas.
# - Simplify the cell styling on the "GWP" sheet (the fill-only style
#   variants are no longer used: header goes back to the plain bold style,
#   body cells go back to the default/unstyled look).
# - Leave the sheet selection/activation state matching the saved file
#   (the "info" sheet ends up active with C14 selected; "GWP" keeps A14
#   selected but is no longer the active tab).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("info")
$ws2 = $wb.Worksheets.Item("GWP")

# --- Restyle the "GWP" sheet: drop the now-unused fill-only xf variants ---
# Header row (was bold-12 + fill) -> back to plain bold-12.
$ws2.Range("A1:G1").Font.Bold = $true
$ws2.Range("A1:G1").Font.Size = 12
# Body rows (was default-font + fill) -> back to the plain default style.
$ws2.Range("A2:G13").Style = "Normal"

# --- "info" sheet: new row 14 for the PE_stream StreamImpactItem ---
$ws1.Range("A14").Value = "PE_stream"
$ws1.Range("B14").Formula = "=B3"
$ws1.Range("C14").Value = "StreamImpactItem"

# --- "GWP" sheet: new row 14 mirroring the PE impact item (row 3) ---
$ws2.Range("A14").Value = "PE_stream"
$ws2.Range("B14").Value = "kg CO2-eq"
$ws2.Range("C14").Formula = "=C3"
$ws2.Range("D14").Formula = "=D3"
$ws2.Range("E14").Formula = "=E3"
$ws2.Range("F14").Formula = "=F3"
$ws2.Range("G14").Formula = "=G3"

# --- Selection / active sheet state ---
$ws2.Activate() | Out-Null
$ws2.Range("A14").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("C14").Select() | Out-Null
